$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2021-11-17"

# Update the row label for the November partial-month row.
$ws.Range("A13").Value = "November (through 11-17)"

# --- Row 13 (November, through 11-17) updated counts / rates ---
# 2016
$ws.Range("F13").Value = 38
$ws.Range("G13").Value = 0.0732
# 2017
$ws.Range("I13").Value = 69
$ws.Range("J13").Value = 0.0282
# 2018
$ws.Range("K13").Value = 6
$ws.Range("L13").Value = 31
$ws.Range("M13").Value = 0.1622
# 2019
$ws.Range("O13").Value = 23
$ws.Range("P13").Value = 0.1481
# 2020
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 96
$ws.Range("S13").Value = 0.04
# 2021
$ws.Range("U13").Value = 118
$ws.Range("V13").Value = 0.008399999999999999

# --- Row 14 (Total) updated counts / rates ---
# 2016
$ws.Range("F14").Value = 472
$ws.Range("G14").Value = 0.1044
# 2017
$ws.Range("I14").Value = 718
$ws.Range("J14").Value = 0.08069999999999999
# 2018
$ws.Range("K14").Value = 72
$ws.Range("L14").Value = 580
$ws.Range("M14").Value = 0.1104
# 2019
$ws.Range("O14").Value = 457
$ws.Range("P14").Value = 0.1022
# 2020
$ws.Range("Q14").Value = 58
$ws.Range("R14").Value = 1099
$ws.Range("S14").Value = 0.0501
# 2021
$ws.Range("U14").Value = 1472
$ws.Range("V14").Value = 0.057
